$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.147.86"
$ws.Range("E2").Value = "  +2.36%  "

$ws.Range("D3").Value = "2.315.92"
$ws.Range("E3").Value = "  +2.13%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.505"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.23%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.18%  "

$ws.Range("E12").Value = "  +3.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +14.49%  "

$ws.Range("E14").Value = "  +3.72%  "

$ws.Range("D15").Value = "2.675.42"
$ws.Range("E15").Value = "  +2.01%  "

$ws.Range("D16").Value = "2.315.49"
$ws.Range("E16").Value = "  +1.95%  "

$ws.Range("E17").Value = "  +4.06%  "

$ws.Range("D18").Value = "43.061.04"
$ws.Range("E18").Value = "  +2.17%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.12%  "

$ws.Range("E20").Value = "  +3.69%  "

$ws.Range("D21").Value = "0.0₃0905"
$ws.Range("E21").Value = "  +1.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.07%  "

$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.51%  "

$ws.Range("E28").Value = "  +7.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.97%  "

$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.24%  "

$ws.Range("E36").Value = "  +3.64%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0693"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "

$ws.Range("E38").Value = "  +4.64%  "

$ws.Range("E39").Value = "  +1.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.80%  "

$ws.Range("E41").Value = "  +1.52%  "

$ws.Range("E42").Value = "  -4.16%  "

$ws.Range("D43").Value = "1.991.43"
$ws.Range("E43").Value = "  +1.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0288"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.17%  "

$ws.Range("E46").Value = "  +1.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.32%  "

$ws.Range("D49").Value = "2.542.49"
$ws.Range("E49").Value = "  +1.93%  "

$ws.Range("E50").Value = "  +4.28%  "

$ws.Range("E51").Value = "  +2.70%  "
